$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume columns in this sheet are plain text in the source data
# (e.g. "219.60", "36.437.47", "  -2.00%  "), not numbers or percentages.
# Column D values look numeric, so Excel would silently reinterpret them as
# numbers (dropping significant trailing zeros, e.g. "219.60" -> 219.6) unless
# the cell is explicitly formatted as Text first. Column E values already have
# surrounding spaces so they are immune to that coercion.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.437.47'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.030.40'
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.55'
$ws.Range("E5").Value = '  -12.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.598'
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '54.59'
$ws.Range("E8").Value = '  -2.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.368'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '56.64'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0749'
$ws.Range("E11").Value = '  -1.86%  '
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.333.08'
$ws.Range("E13").Value = '  +2.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.21'
$ws.Range("E14").Value = '  -0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '19.84'
$ws.Range("E15").Value = '  -9.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.754'
$ws.Range("E16").Value = '  -2.51%  '
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.018.81'
$ws.Range("E18").Value = '  +1.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.377.81'
$ws.Range("E19").Value = '  -1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.78'
$ws.Range("E20").Value = '  +13.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.23'
$ws.Range("E21").Value = '  -3.70%  '
$ws.Range("E22").Value = '  -3.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '219.60'
$ws.Range("E23").Value = '  -6.02%  '
$ws.Range("E24").Value = '  -0.12%  '
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  -8.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.18'
$ws.Range("E27").Value = '  -1.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.66'
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("E29").Value = '  -1.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '18.83'
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("E31").Value = '  +2.61%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.33'
$ws.Range("E33").Value = '  -4.79%  '
$ws.Range("E34").Value = '  -3.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.46'
$ws.Range("E35").Value = '  +1.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.23'
$ws.Range("E36").Value = '  -3.16%  '
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -2.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.74'
$ws.Range("E39").Value = '  +6.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.23'
$ws.Range("E40").Value = '  -8.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.52'
$ws.Range("E41").Value = '  +44.91%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.94'
$ws.Range("E42").Value = '  -3.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.484.23'
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0932'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.63'
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0202'
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("E47").Value = '  -5.98%  '
$ws.Range("E48").Value = '  -0.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.997'
$ws.Range("E49").Value = '  -2.98%  '
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.91'
$ws.Range("E51").Value = '  +1.29%  '
